# Adds the three "Formatted <category>" summary sheets (services, technologies,
# sectors) produced by the first phase of text analysis. Each sheet gets:
#   - B1: "Formatted <category>" header, C1: "count" header (bold, bordered,
#     centered/top-aligned — matching the style already used by the existing
#     "before"/"cleaned"/"Cleaned_FULL" sheets' header + index cells)
#   - column A: a 0-based row index, same style as the headers
#   - column B: the formatted label (plain)
#   - column C: the count (plain, numeric)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

function Add-DataSheet($Name, $HeaderLabel, $Rows) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $Name

    $n = $Rows.Count
    $lastRow = $n + 1

    # Clone the existing header/index style (bold font, thin border all
    # round, centered horizontally + top vertically) from sheet 1's A2 cell
    # rather than re-deriving the font/border/alignment by hand, so the
    # new cells land on the very same cellXfs entry the workbook already
    # uses for this look.
    $ws1.Range("A2").Copy()
    $ws.Range("B1:C1").PasteSpecial(-4122)
    $ws.Range("A2:A$lastRow").PasteSpecial(-4122)

    $ws.Cells.Item(1,2).Value = $HeaderLabel
    $ws.Cells.Item(1,3).Value = "count"

    for ($i = 0; $i -lt $n; $i++) {
        $r = $i + 2
        $label = $Rows[$i][0]
        $count = $Rows[$i][1]
        $ws.Cells.Item($r,1).Value = $i
        $ws.Cells.Item($r,2).Value = $label
        $ws.Cells.Item($r,3).Value = $count
    }
}

$data_services = @(
    @('sme support', 247),
    @('technological innovation', 216),
    @('knowledge transfer', 210),
    @('ecosystem building', 203),
    @('technology transfer', 201),
    @('innovation management', 199),
    @('prototyping', 189),
    @('sme business development', 159),
    @('public sector innovation', 144),
    @('finance', 132),
    @('regional development', 129),
    @('circular economy', 119),
    @('vocational training', 118),
    @('smart specialisation strategies', 105),
    @('field trial', 99)
)
Add-DataSheet "services" "Formatted services" $data_services

$data_technologies = @(
    @('artificial intelligence & decision support', 268),
    @('internet of things', 207),
    @('cybersecurity', 204),
    @('big data', 173),
    @('digital twins', 166),
    @('robotics', 158),
    @('virtual reality', 146),
    @('high performance computing', 146),
    @('sensors & vision processing systems', 131),
    @('simulation engineering and modelling', 118),
    @('additive manufacturing', 112),
    @('cloud services', 103),
    @('cyber-physical systems', 102),
    @('internet services & applications', 92),
    @('communication network', 90),
    @('blockchain and distributed ledger technology (dlt)', 89),
    @('human computer interaction', 84),
    @('mobility', 79),
    @('software architectures', 77),
    @('logistics', 66),
    @('location-based applications', 57),
    @('laser-based manufacturing and materials processing', 54),
    @('gamification', 54),
    @('quantum technologies (computing/communication)', 47),
    @('optoelectronics', 41),
    @('micro- and nanoelectronics', 41),
    @('semiconductors and nanotechnology', 36),
    @('new technologies for audio-visual sector - media', 34),
    @('industrial biotechnology', 31),
    @('bi tools', 28),
    @('photonics', 27),
    @('products)', 24),
    @('chemical engineering (plants', 24),
    @('displays', 17),
    @('organic and large area electronics', 14),
    @('data', 14)
)
Add-DataSheet "technologies" "Formatted technologies" $data_technologies

$data_sectors = @(
    @('manufacturing and processing', 196),
    @('energy', 152),
    @('health care', 140),
    @('smart city', 134),
    @('public administration', 128),
    @('education', 123),
    @('environment', 117),
    @('transport & mobility', 115),
    @('automotive', 109),
    @('agricultural biotechnology and food biotechnology', 95),
    @('construction & assembly', 91),
    @('metal working and industrial production', 87),
    @('travel and tourism', 78),
    @('telecommunications', 74),
    @('food and beverages', 66),
    @('life sciences', 65),
    @('security', 64),
    @('cultural and creative economy', 62),
    @('retail', 56),
    @('wholesale or distribution', 56),
    @('financial', 55),
    @('community-led local development', 54),
    @('textiles', 51),
    @('transport sector', 50),
    @('consumer products', 47),
    @('personal services', 45),
    @('maritime', 44),
    @('defence', 43),
    @('polymers and plastics', 42),
    @('aeronautics', 42),
    @('space', 41),
    @('electricity', 33),
    @('nmp non-metallic materials & basic processes', 29),
    @('paper and wood', 28),
    @('fuels and petroleum engineering', 26),
    @('real estate', 21),
    @('legal aspects', 16),
    @('regulation', 15),
    @('fishery', 15),
    @('mining and extraction', 12),
    @('leather', 11),
    @('nuclear', 8),
    @('tobacco', 4)
)
Add-DataSheet "sectors" "Formatted sectors" $data_sectors
